# Weekly fruit/vegetable price update: insert a new daily record at the top
# of the "Mango" dataset (row 178), shifting all subsequent records down by
# one row (old row 178 -> new row 179, ..., old row 298 -> new row 299).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 178; this shifts rows 178:298 down to
# 179:299 and extends the sheet dimension to A1:T299 automatically.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with the new record's values.
$ws.Cells.Item(178, 1).Value = 10
$ws.Cells.Item(178, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(178, 3).Value = "La Araucanía"
$ws.Cells.Item(178, 4).Value = 44603
$ws.Cells.Item(178, 5).Value = 9
$ws.Cells.Item(178, 6).Value = "Fruta"
$ws.Cells.Item(178, 7).Value = 100108
$ws.Cells.Item(178, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(178, 9).Value = 100108002
$ws.Cells.Item(178, 10).Value = "Mango"
$ws.Cells.Item(178, 11).Value = "Sin especificar"
$ws.Cells.Item(178, 12).Value = "Primera"
$ws.Cells.Item(178, 13).Value = 300
$ws.Cells.Item(178, 14).Value = 8000
$ws.Cells.Item(178, 15).Value = 8000
$ws.Cells.Item(178, 16).Value = 8000
$ws.Cells.Item(178, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(178, 18).Value = "Perú"
$ws.Cells.Item(178, 19).Value = 2000
$ws.Cells.Item(178, 20).Value = 4
